# Egypt COVID-19 data update (commit: "data update jun 13th")
# Extends the daily case/recovery/death series through 2020-06-13
# (serial date 43995) and pads the sheet with blank formatted rows
# beyond the data, matching the author's habitual layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New daily rows 82:94 (2020-06-01 .. 2020-06-13) didn't exist before.
#    Seed them by copying row 81's cell (date format m/d/yyyy, style s="1")
#    down through row 111 - this also pre-creates the blank pad rows
#    95:111 with the same date-column formatting/style as the data rows.
# ---------------------------------------------------------------------
$ws.Range("A81").Copy($ws.Range("A82:A111"))

# ---------------------------------------------------------------------
# 2. Write the new date serials into A82:A94 (rows 95:111 stay blank -
#    they only keep the number-format/style copied above).
# ---------------------------------------------------------------------
$dates = @(
    43983,
    43984,
    43985,
    43986,
    43987,
    43988,
    43989,
    43990,
    43991,
    43992,
    43993,
    43994,
    43995
)
$arrA = New-Object 'object[,]' $dates.Count,1
for ($i = 0; $i -lt $dates.Count; $i++) {
    $arrA[$i,0] = $dates[$i]
}
$ws.Range("A82:A94").Value = $arrA

# Drop the date values that the copy-down put into the trailing pad
# rows - only the style/number-format should survive there.
$ws.Range("A95:A111").ClearContents()

# ---------------------------------------------------------------------
# 3. Fill in new_cases / all_cases / recovered / all_deaths / new_deaths
#    / converters (columns B:G) for every day from 2020-03-15 (row 64,
#    previously date-only) through 2020-06-13 (row 94, brand new).
# ---------------------------------------------------------------------
$rows = @(
    @(398,10829,2626,571,15,3133),
    @(399,11228,2799,592,21,3363),
    @(491,11719,2950,612,20,3526),
    @(510,12229,3172,630,18,3742),
    @(535,12764,3440,645,15,4001),
    @(720,13484,3742,659,14,4275),
    @(745,14229,3994,680,21,4584),
    @(774,15003,4217,696,16,4798),
    @(783,15786,4374,707,11,4960),
    @(727,16513,4628,735,28,5192),
    @(752,17265,4807,764,29,5366),
    @(702,17967,4900,783,19,5481),
    @(789,18756,5027,797,14,5606),
    @(910,19666,5205,816,19,5798),
    @(1127,20793,5359,845,29,6019),
    @(1289,22082,5511,879,34,6237),
    @(1367,23449,5693,913,34,6456),
    @(1536,24985,6037,959,46,6810),
    @(1399,26384,6447,1005,46,7149),
    @(1152,27536,6827,1052,47,7642),
    @(1079,28615,7350,1088,36,8371),
    @(1152,29767,7756,1126,38,8793),
    @(1348,31115,8158,1166,40,9216),
    @(1497,32612,8538,1198,32,9603),
    @(1467,34079,8961,1237,39,10131),
    @(1365,35444,9375,1271,34,10618),
    @(1385,36829,9786,1306,35,11071),
    @(1455,38284,10289,1342,36,11583),
    @(1442,39726,10691,1377,35,12062),
    @(1577,41303,11108,1422,45,12493),
    @(1677,42980,11529,1484,62,12919)
)
$arrBG = New-Object 'object[,]' $rows.Count,6
for ($i = 0; $i -lt $rows.Count; $i++) {
    for ($j = 0; $j -lt 6; $j++) {
        $arrBG[$i,$j] = $rows[$i][$j]
    }
}
$ws.Range("B64:G94").Value = $arrBG

# ---------------------------------------------------------------------
# 4. Restore the author's view state: scrolled down to the new data and
#    with G94 (the last new cell touched) selected.
# ---------------------------------------------------------------------
$ws.Range("G94").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 92
$win.ScrollColumn = 1
